$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case block "B.8" (rows 47-50), appended below the existing "B.7"
# block, following the same layout used throughout the sheet:
#   col B -> Test ID            (centered)
#   col C -> Test Description
#   col D -> Steps
#   col E -> Expected Results
#   col F -> Status

$ws.Range("B47").Value = "B.8"
$ws.Range("C47").Value = "Verify if the inventory gets refreshed on adding and deleting "
$ws.Range("D47").Value = "1.Launch the URL"

$ws.Range("C48").Value = "an item"
$ws.Range("D48").Value = "2.Login to the application"
$ws.Range("E48").Value = "The inventory should get refreshed on adding and deleting "
$ws.Range("F48").Value = "In progress"

$ws.Range("D49").Value = "3.Verify if the inventory gets refreshed on adding "
$ws.Range("E49").Value = "items"

# D50 stays blank (end of block spacer row), same as D46/D41 above it.

# Pick up the same formatting (font/alignment) already used by the matching
# columns of the preceding test-case rows, instead of leaving the new cells
# on the workbook default style.
$ws.Range("B37").Copy()
$ws.Range("B47").PasteSpecial(-4122)

$ws.Range("C37").Copy()
$ws.Range("C47").PasteSpecial(-4122)

$ws.Range("D36").Copy()
$ws.Range("D47").PasteSpecial(-4122)

$ws.Range("C38").Copy()
$ws.Range("C48").PasteSpecial(-4122)

$ws.Range("D37").Copy()
$ws.Range("D48").PasteSpecial(-4122)

$ws.Range("E36").Copy()
$ws.Range("E48").PasteSpecial(-4122)

$ws.Range("F36").Copy()
$ws.Range("F48").PasteSpecial(-4122)

$ws.Range("D39").Copy()
$ws.Range("D49").PasteSpecial(-4122)

$ws.Range("E37").Copy()
$ws.Range("E49").PasteSpecial(-4122)

$ws.Range("D40").Copy()
$ws.Range("D50").PasteSpecial(-4122)

# Center the new Test ID cell, matching B16/B22/B27/B31/B37/B42.
$ws.Range("B47").HorizontalAlignment = -4108
